{"js": "// Move the \"_GoBack\" bookmark from the end of the \"Suburb\" merge-field\n// paragraph to the end of the \"To\" paragraph, leaving all other content\n// (the MERGEFIELD runs for First_Name/Last_Name, Street_Address, Suburb,\n// etc.) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"To\" paragraph and the \"Suburb\" merge-field paragraph by\n// their (loaded) plain text content, so the script is resilient to any\n// paragraphs that might shift position.\nlet toParaIndex = -1;\nlet suburbParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text.trim();\n  if (toParaIndex === -1 && t === \"To\") {\n    toParaIndex = i;\n  }\n  if (t.indexOf(\"\\u00ABSuburb\\u00BB\") !== -1) {\n    suburbParaIndex = i;\n  }\n}\n\nif (toParaIndex === -1 || suburbParaIndex === -1) {\n  throw new Error(\"Could not locate the 'To' and/or 'Suburb' paragraphs.\");\n}\n\n// Remove the existing bookmark (currently sitting at the end of the\n// Suburb paragraph) \u2026\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// \u2026 and re-insert it, collapsed, at the end of the \"To\" paragraph.\nconst toEnd = paragraphs.items[toParaIndex].getRange(\"End\");\ntoEnd.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from the end of the \"Suburb\" merge-field\n# paragraph to the end of the \"To\" paragraph, leaving all other content\n# untouched.\n\n$d = $word.ActiveDocument\n\n# Remove the bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Find the \"To\" paragraph and collapse a range to its end (just before the\n# paragraph mark), then drop the bookmark there.\n$toParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"To\") {\n        $toParagraph = $p\n        break\n    }\n}\n\n$target = $toParagraph.Range\n$target.SetRange($target.Start, $target.End - 1)\n$target.Collapse(0)  # wdCollapseEnd\n\n$d.Bookmarks.Add(\"_GoBack\", $target)\n"}
